# "filtro + modifica generatore excel"
# The export now filters out disabled/USER-role rows, so the "Leo" (id=1,
# USER) row disappears and every following row shifts up one. The export
# generator was also tweaked: a freshly registered user ("Fabrizio") is
# appended at the bottom with a blank password column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the filtered-out row (id=1, Leo, USER) - remaining rows shift up.
$ws.Rows(2).Delete()

# Re-assert the values for the rows now in their final positions (the
# regenerated export carries slightly different id/created_at stamps).
$ws.Range("A2").Value = 2
$ws.Range("B2").Value = "Edo"
$ws.Range("C2").Value = 45061.763309236114
$ws.Range("D2").Value = "rich"
$ws.Range("E2").Value = "Edo"
$ws.Range("F2").Value = "ADMIN"

$ws.Range("A3").Value = 3
$ws.Range("B3").Value = "Ettore"
$ws.Range("C3").Value = 45069.87578456019
$ws.Range("D3").Value = "rich"
$ws.Range("E3").Value = "Ettore"
$ws.Range("F3").Value = "ADMIN"

$ws.Range("A4").Value = 4
$ws.Range("B4").Value = "Giacomo"
$ws.Range("C4").Value = 45070.984361331015
$ws.Range("D4").Value = "rich"
$ws.Range("E4").Value = "Giacomo"
$ws.Range("F4").Value = "ADMIN"

# New user appended by the regenerated export - no password set yet.
$ws.Range("A5").Value = 6
$ws.Range("B5").Value = "Fabrizio"
$ws.Range("C5").Value = 45070.9851908912
$ws.Range("D5").Value = "rich"
$ws.Range("E5").Value = ""
$ws.Range("F5").Value = "ADMIN"

# Match the "created_at" date formatting (column C) on the newly written rows.
$ws.Range("C2").Copy()
$ws.Range("C4:C5").PasteSpecial(-4122)
$excel.CutCopyMode = $false
